$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B43: change from text "3" to numeric 3
$ws.Range("B43").Value = 3

# Add new row 44 with data from the diff
$ws.Range("A44").Value = "Ruilin"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "3"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "无"
$ws.Range("D44").Value = "CRT"
$ws.Range("E44").Value = "OTH"
$ws.Range("F44").Value = "17635cfa-5d3f-4715-99a5-e710de1fbea7"
$ws.Range("G44").Value = "S1XXq6lRW_annotated.xlsx"
$ws.Range("H44").Value = "Technical contribution of the paper is very limited."
